$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '[Jingya%Wang%NULL%1,                           Jennifer M%Cooper%NULL%2,                           Jennifer M%Cooper%NULL%0,                           Krishna%Gokhale%NULL%1,                           Dionisio%Acosta-Mena%NULL%1,                           Samir%Dhalla%NULL%1,                           Nathan%Byne%NULL%1,                           Joht Singh%Chandan%NULL%1,                           Astha%Anand%NULL%2,                           Astha%Anand%NULL%0,                           Kelvin%Okoth%NULL%1,                           Anuradhaa%Subramanian%NULL%2,                           Anuradhaa%Subramanian%NULL%0,                           Mansoor N%Bangash%NULL%1,                           Thomas%Jackson%NULL%2,                           Thomas%Jackson%NULL%0,                           Dawit%Zemedikun%NULL%2,                           Dawit%Zemedikun%NULL%0,                           Tom%Taverner%NULL%2,                           Tom%Taverner%NULL%0,                           Wasim%Hanif%NULL%1,                           Sandip%Ghosh%NULL%1,                           Parth%Narendran%NULL%1,                           Konstantinos A%Toulis%NULL%1,                           Abd A%Tahrani%NULL%2,                           Abd A%Tahrani%NULL%0,                           Rajendra%Surenthirakumaran%NULL%2,                           Rajendra%Surenthirakumaran%NULL%0,                           Nicola J%Adderley%NULL%1,                           Shamil%Haroon%NULL%2,                           Shamil%Haroon%NULL%0,                           Kamlesh%Khunti%NULL%4,                           Kamlesh%Khunti%NULL%0,                           Christopher%Sainsbury%NULL%2,                           Christopher%Sainsbury%NULL%0,                           G Neil%Thomas%NULL%1,                           Krishnarajah%Nirantharakumar%NULL%2,                           Krishnarajah%Nirantharakumar%NULL%0]'
$ws.Range("E3").Value = '[Laura%Orioli%NULL%1,                           Thomas%Servais%NULL%1,                           Leïla%Belkhir%NULL%1,                           Pierre-François%Laterre%NULL%1,                           Jean-Paul%Thissen%NULL%1,                           Bernard%Vandeleene%NULL%1,                           Dominique%Maiter%NULL%1,                           Jean C.%Yombi%NULL%1,                           Michel P.%Hermans%NULL%1]'
$ws.Range("E4").Value = '[Jean-Daniel%Lalau%NULL%1,                           Abdallah%Al-Salameh%NULL%3,                           Samy%Hadjadj%NULL%3,                           Thomas%Goronflot%NULL%1,                           Nicolas%Wiernsperger%NULL%1,                           Matthieu%Pichelin%NULL%3,                           Ingrid%Allix%NULL%3,                           Coralie%Amadou%NULL%2,                           Olivier%Bourron%NULL%2,                           Thierry%Duriez%NULL%1,                           Jean-François%Gautier%NULL%3,                           Anne%Dutour%NULL%1,                           Céline%Gonfroy%NULL%1,                           Didier%Gouet%NULL%1,                           Michael%Joubert%NULL%3,                           Ingrid%Julier%NULL%1,                           Etienne%Larger%NULL%1,                           Lucien%Marchand%NULL%3,                           Michel%Marre%NULL%1,                           Laurent%Meyer%NULL%3,                           Frédérique%Olivier%NULL%1,                           Gaëtan%Prevost%NULL%3,                           Pascale%Quiniou%NULL%1,                           Christelle%Raffaitin-Cardin%NULL%1,                           Ronan%Roussel%NULL%3,                           Pierre-Jean%Saulnier%NULL%3,                           Dominique%Seret-Begue%NULL%1,                           Charles%Thivolet%NULL%3,                           Camille%Vatier%NULL%2,                           Rachel%Desailloud%NULL%1,                           Matthieu%Wargny%NULL%3,                           Pierre%Gourdy%NULL%3,                           Bertrand%Cariou%NULL%1,                           NULL%NULL%NULL%39]'
$ws.Range("E5").Value = '[Nan%Jiang%NULL%1,                           Zhenyuan%Chen%NULL%1,                           Li%Liu%NULL%3,                           Xiaoxv%Yin%NULL%1,                           Heping%Yang%NULL%1,                           Xiangping%Tan%NULL%1,                           Jing%Wang%NULL%2,                           Hui%Li%NULL%13,                           Mengge%Tian%NULL%1,                           Zuxun%Lu%NULL%1,                           Nian%Xiong%NULL%4,                           Yanhong%Gong%NULL%1]'
$ws.Range("E6").Value = '[Xu%Cheng%NULL%1,                           Ye-Mao%Liu%NULL%2,                           Haomiao%Li%NULL%1,                           Xin%Zhang%NULL%1,                           Fang%Lei%NULL%2,                           Juan-Juan%Qin%NULL%2,                           Ze%Chen%NULL%1,                           Ke-Qiong%Deng%NULL%1,                           Lijin%Lin%NULL%2,                           Ming-Ming%Chen%NULL%2,                           Xiaohui%Song%NULL%1,                           Meng%Xia%NULL%2,                           Xuewei%Huang%NULL%2,                           Weifang%Liu%NULL%1,                           Jingjing%Cai%NULL%2,                           Xiao-Jing%Zhang%NULL%2,                           Feng%Zhou%NULL%1,                           Peng%Zhang%NULL%1,                           Yibin%Wang%NULL%2,                           Xinliang%Ma%NULL%1,                           Qingbo%Xu%NULL%2,                           Juan%Yang%NULL%3,                           Ping%Ye%NULL%2,                           Weiming%Mao%NULL%2,                           Xiaodong%Huang%NULL%2,                           Jiahong%Xia%NULL%2,                           Bing-Hong%Zhang%NULL%2,                           Jiao%Guo%NULL%1,                           Lihua%Zhu%NULL%2,                           Zhibing%Lu%NULL%3,                           Yufeng%Yuan%NULL%2,                           Xiang%Wei%NULL%1,                           Zhi-Gang%She%NULL%2,                           Yan-Xiao%Ji%NULL%2,                           Hongliang%Li%NULL%2]'
$ws.Range("E7").Value = '[Zhelong%Liu%NULL%1,                           Xi%Bai%NULL%1,                           Xia%Han%NULL%1,                           Wangyan%Jiang%NULL%1,                           Lin%Qiu%NULL%2,                           Shi%Chen%NULL%1,                           Xuefeng%Yu%NULL%1]'
$ws.Range("E8").Value = '[Qianhui%Zhang%NULL%1,                           Yanhong%Wei%NULL%1,                           Min%Chen%NULL%1,                           Qianqian%Wan%NULL%1,                           Xiaoqi%Chen%NULL%1]'
$ws.Range("E9").Value = '[G. Antonio%Silverii%NULL%1,                           Matteo%Monami%NULL%1,                           Achille%Cernigliaro%NULL%1,                           Enrica%Vigneri%NULL%1,                           Valentina%Guarnotta%NULL%1,                           Salvatore%Scondotto%NULL%1,                           Vincenza A.%Allotta%NULL%1,                           Michela%Conti%NULL%1,                           Carla%Giordano%NULL%1,                           Edoardo%Mannucci%NULL%1]'
$ws.Range("E10").Value = '[Jiaoyue%Zhang%NULL%1,                           Wen%Kong%NULL%1,                           Pengfei%Xia%NULL%1,                           Ying%Xu%NULL%1,                           Li%Li%NULL%3,                           Qin%Li%NULL%1,                           Li%Yang%NULL%3,                           Qi%Wei%NULL%3,                           Hanyu%Wang%NULL%1,                           Huiqing%Li%NULL%1,                           Juan%Zheng%NULL%1,                           Hui%Sun%NULL%1,                           Wenfang%Xia%NULL%1,                           Geng%Liu%NULL%1,                           Xueyu%Zhong%NULL%1,                           Kangli%Qiu%NULL%1,                           Yan%Li%NULL%1,                           Han%Wang%NULL%1,                           Yuxiu%Wang%NULL%1,                           Xiaoli%Song%NULL%1,                           Hua%Liu%NULL%1,                           Si%Xiong%NULL%1,                           Yumei%Liu%NULL%1,                           Zhenhai%Cui%NULL%1,                           Yu%Hu%NULL%6,                           Lulu%Chen%NULL%1,                           An%Pan%NULL%1,                           Tianshu%Zeng%NULL%1]'
$ws.Range("E11").Value = '[Jinghong%Li%NULL%2,                           Qi%Wei%NULL%0,                           Willis X%Li%NULL%1,                           Karen C%McCowen%NULL%1,                           Wei%Xiong%NULL%2,                           Jiao%Liu%NULL%2,                           Wenlijun%Jiang%NULL%2,                           Traci%Marin%NULL%1,                           Robert L%Thomas%NULL%1,                           Ming%He%NULL%2,                           Brendan%Gongol%NULL%1,                           Mark%Hepokoski%NULL%2,                           Jason X-J%Yuan%NULL%1,                           John Y-J%Shyy%NULL%1,                           Nian%Xiong%NULL%0,                           Atul%Malhotra%NULL%2]'
$ws.Range("E12").Value = '[Xueqi%Cheng%NULL%1,                           Siyi%Xin%NULL%1,                           Yaqi%Chen%NULL%1,                           Leyu%Li%NULL%1,                           Wanjun%Chen%NULL%1,                           Wenjia%Li%NULL%1,                           Baoan%Zhou%NULL%1,                           Chenxia%Li%NULL%1,                           Yu%Gong%NULL%1,                           Fei%Li%NULL%1,                           Peng%Duan%NULL%1,                           Xingjian%Zhou%NULL%1]'
$ws.Range("E13").Value = '[Carolyn T.%Bramante%bramante@umn.edu%0,                           John%Buse%NULL%2,                           John%Buse%NULL%0,                           Leonardo%Tamaritz%NULL%1,                           Ana%Palacio%NULL%2,                           Ken%Cohen%NULL%1,                           Deneen%Vojta%NULL%2,                           David%Liebovitz%NULL%1,                           Nia%Mitchell%NULL%1,                           Jacinda%Nicklas%NULL%1,                           Ildiko%Lingvay%NULL%1,                           Jeanne M.%Clark%NULL%1,                           Louis J.%Aronne%NULL%1,                           Erik%Anderson%NULL%1,                           Michael%Usher%NULL%1,                           Ryan%Demmer%NULL%1,                           Genevieve B.%Melton%NULL%1,                           Nicholas%Ingraham%NULL%1,                           Christopher J.%Tignanelli%NULL%2,                           Christopher J.%Tignanelli%NULL%0]'
$ws.Range("E14").Value = '[Tak Kyu%Oh%NULL%2,                           In-Ae%Song%songoficu@outlook.kr%1]'
$ws.Range("E15").Value = '[Reyan%Ghany%NULL%1,                           Ana%Palacio%NULL%0,                           Elissa%Dawkins%NULL%1,                           Gordon%Chen%NULL%1,                           Daniel%McCarter%NULL%1,                           Emancia%Forbes%NULL%1,                           Brian%Chung%NULL%1,                           Leonardo%Tamariz%NULL%2]'
$ws.Range("E16").Value = '[Jinghong%Li%NULL%0,                           Qi%Wei%NULL%0,                           Karen C.%McCowen%NULL%1,                           Wei%Xiong%NULL%0,                           Jiao%Liu%NULL%0,                           Wenlijun%Jiang%NULL%0,                           Robert L.%Thomas%NULL%1,                           Mark%Hepokoski%NULL%0,                           Ming%He%NULL%0,                           John Y. J.%Shyy%NULL%1,                           Atul%Malhotra%NULL%0,                           Nian%Xiong%nianxiong@hust.edu.cn%0,                           Willis X.%Li%wxli@health.ucsd.edu%1]'
$ws.Range("E17").Value = '[J. W.%Goodall%NULL%1,                           T. A. N.%Reed%NULL%2,                           T. A. N.%Reed%NULL%0,                           M.%Ardissino%NULL%1,                           P.%Bassett%NULL%1,                           A. M.%Whittington%NULL%1,                           D. L.%Cohen%NULL%1,                           N.%Vaid%NULL%1]'
$ws.Range("E18").Value = '[Chioma%Izzi-Engbeaya%NULL%1,                           Walter%Distaso%NULL%2,                           Walter%Distaso%NULL%0,                           Anjali%Amin%NULL%1,                           Wei%Yang%NULL%1,                           Oluwagbemiga%Idowu%NULL%1,                           Julia S%Kenkre%NULL%1,                           Ronak J%Shah%NULL%1,                           Evelina%Woin%NULL%1,                           Christine%Shi%NULL%1,                           Nael%Alavi%NULL%2,                           Nael%Alavi%NULL%0,                           Hala%Bedri%NULL%1,                           Niamh%Brady%NULL%1,                           Sophie%Blackburn%NULL%1,                           Martina%Leczycka%NULL%1,                           Sanya%Patel%NULL%1,                           Elizaveta%Sokol%NULL%1,                           Edward%Toke-Bjolgerud%NULL%1,                           Ambreen%Qayum%NULL%1,                           Mariana%Abdel-Malek%NULL%1,                           David C D%Hope%NULL%1,                           Nick S%Oliver%NULL%1,                           Vasiliki%Bravis%NULL%1,                           Shivani%Misra%NULL%1,                           Tricia M%Tan%NULL%1,                           Neil E%Hill%NULL%2,                           Neil E%Hill%NULL%0,                           Victoria%Salem%NULL%2,                           Victoria%Salem%NULL%0]'
$ws.Range("E19").Value = '[So-Yeon%Kim%NULL%1,                           Kyung-Soo%Kim%NULL%1]'
$ws.Range("E20").Value = '[Andrew B.%Crouse%NULL%2,                           Tiffany%Grimes%NULL%1,                           Peng%Li%NULL%1,                           Matthew%Might%NULL%1,                           Fernando%Ovalle%NULL%1,                           Anath%Shalev%NULL%1]'
$ws.Range("E21").Value = '[Carolyn T%Bramante%NULL%1,                           Nicholas E%Ingraham%NULL%1,                           Thomas A%Murray%NULL%1,                           Schelomo%Marmor%NULL%1,                           Shane%Hovertsen%NULL%1,                           Jessica%Gronski%NULL%1,                           Chace%McNeil%NULL%1,                           Ruoying%Feng%NULL%1,                           Gabriel%Guzman%NULL%1,                           Nermine%Abdelwahab%NULL%1,                           Samantha%King%NULL%1,                           Leonardo%Tamariz%NULL%0,                           Thomas%Meehan%NULL%1,                           Kathryn M%Pendleton%NULL%1,                           Bradley%Benson%NULL%1,                           Deneen%Vojta%NULL%0,                           Christopher J%Tignanelli%NULL%1]'
$ws.Range("E22").Value = '[Bertrand%Cariou%bertrand.cariou@univ-nantes.fr%1,                           Samy%Hadjadj%samy.hadjadj@univ-nantes.fr%0,                           Samy%Hadjadj%samy.hadjadj@univ-nantes.fr%0,                           Matthieu%Wargny%NULL%0,                           Matthieu%Wargny%NULL%0,                           Matthieu%Pichelin%NULL%0,                           Matthieu%Pichelin%NULL%0,                           Abdallah%Al-Salameh%NULL%0,                           Abdallah%Al-Salameh%NULL%0,                           Ingrid%Allix%NULL%0,                           Ingrid%Allix%NULL%0,                           Coralie%Amadou%NULL%0,                           Gwénaëlle%Arnault%NULL%2,                           Gwénaëlle%Arnault%NULL%0,                           Florence%Baudoux%NULL%1,                           Bernard%Bauduceau%NULL%1,                           Sophie%Borot%NULL%1,                           Muriel%Bourgeon-Ghittori%NULL%2,                           Muriel%Bourgeon-Ghittori%NULL%0,                           Olivier%Bourron%NULL%0,                           David%Boutoille%NULL%4,                           David%Boutoille%NULL%0,                           France%Cazenave-Roblot%NULL%3,                           France%Cazenave-Roblot%NULL%0,                           Claude%Chaumeil%NULL%2,                           Claude%Chaumeil%NULL%0,                           Emmanuel%Cosson%NULL%1,                           Sandrine%Coudol%NULL%2,                           Sandrine%Coudol%NULL%0,                           Patrice%Darmon%NULL%1,                           Emmanuel%Disse%NULL%2,                           Emmanuel%Disse%NULL%0,                           Amélie%Ducet-Boiffard%NULL%2,                           Amélie%Ducet-Boiffard%NULL%0,                           Bénédicte%Gaborit%NULL%1,                           Michael%Joubert%NULL%0,                           Michael%Joubert%NULL%0,                           Véronique%Kerlan%NULL%2,                           Véronique%Kerlan%NULL%0,                           Bruno%Laviolle%NULL%1,                           Lucien%Marchand%NULL%0,                           Lucien%Marchand%NULL%0,                           Laurent%Meyer%NULL%0,                           Laurent%Meyer%NULL%0,                           Louis%Potier%NULL%2,                           Louis%Potier%NULL%0,                           Gaëtan%Prevost%NULL%0,                           Gaëtan%Prevost%NULL%0,                           Jean-Pierre%Riveline%NULL%1,                           René%Robert%NULL%2,                           René%Robert%NULL%0,                           Pierre-Jean%Saulnier%NULL%0,                           Pierre-Jean%Saulnier%NULL%0,                           Ariane%Sultan%NULL%2,                           Ariane%Sultan%NULL%0,                           Jean-François%Thébaut%NULL%2,                           Jean-François%Thébaut%NULL%0,                           Charles%Thivolet%NULL%0,                           Charles%Thivolet%NULL%0,                           Blandine%Tramunt%NULL%2,                           Blandine%Tramunt%NULL%0,                           Camille%Vatier%NULL%0,                           Ronan%Roussel%NULL%0,                           Ronan%Roussel%NULL%0,                           Jean-François%Gautier%NULL%0,                           Jean-François%Gautier%NULL%0,                           Pierre%Gourdy%NULL%0,                           Pierre%Gourdy%NULL%0,                           NULL%NULL%NULL%0,                           NULL%NULL%NULL%0]'
$ws.Range("E23").Value = '[Kamlesh%Khunti%NULL%1,                           Peter%Knighton%NULL%1,                           Francesco%Zaccardi%NULL%2,                           Chirag%Bakhai%NULL%1,                           Emma%Barron%NULL%1,                           Naomi%Holman%NULL%1,                           Partha%Kar%NULL%1,                           Claire%Meace%NULL%1,                           Naveed%Sattar%NULL%1,                           Stephen%Sharp%NULL%1,                           Nicholas J%Wareham%NULL%1,                           Andy%Weaver%NULL%1,                           Emilia%Woch%NULL%1,                           Bob%Young%NULL%1,                           Jonathan%Valabhji%NULL%2]'
$ws.Range("E24").Value = '[Luis M.%Pérez-Belmonte%luismiguelpb1984@gmail.com%1,                           José David%Torres-Peña%NULL%1,                           María D.%López-Carmona%NULL%1,                           M. Mar.%Ayala-Gutiérrez%NULL%1,                           Francisco%Fuentes-Jiménez%NULL%1,                           Lucía Jorge%Huerta%NULL%1,                           Jaime Alonso%Muñoz%NULL%1,                           Manuel%Rubio-Rivas%NULL%1,                           Manel%Madrazo%NULL%1,                           Marcos Guzmán%Garcia%NULL%1,                           Beatriz Vicente%Montes%NULL%1,                           Joaquim Fernández%Sola%NULL%1,                           Javier%Ena%NULL%1,                           Ruth Gonzalez%Ferrer%NULL%1,                           Carmen Mella%Pérez%NULL%1,                           Carlos Jorge%Ripper%NULL%1,                           Jose Javier Napal%Lecumberri%NULL%1,                           Iris El Attar%Acedo%NULL%1,                           Susana Plaza%Canteli%NULL%1,                           Sara Fuente%Cosío%NULL%1,                           Francisco Amorós%Martínez%NULL%1,                           Begoña Cortés%Rodríguez%NULL%1,                           Pablo%Pérez-Martínez%pablopermar@yahoo.es%1,                           José Manuel%Ramos-Rincón%NULL%1,                           Ricardo%Gómez-Huelgas%NULL%1,                           NULL%NULL%NULL%0]'
$ws.Range("E25").Value = '[Jun Young%Do%NULL%1,                           Sang Won%Kim%NULL%1,                           Jong Won%Park%NULL%1,                           Kyu Hyang%Cho%NULL%1,                           Seok Hui%Kang%NULL%1]'
$ws.Range("E26").Value = '[Yuchen%Chen%xref no email%0,       Dong%Yang%xref no email%0,       Biao%Cheng%xref no email%0,       Jian%Chen%xref no email%0,       Anlin%Peng%xref no email%0,       Chen%Yang%xref no email%0,       Chong%Liu%xref no email%0,       Mingrui%Xiong%xref no email%0,       Aiping%Deng%xref no email%0,       Yu%Zhang%xref no email%0,       Ling%Zheng%xref no email%0,       Kun%Huang%xref no email%0]'
$ws.Range("E27").Value = '[Michelle A.%Lally%NULL%1,                           Philip%Tsoukas%NULL%1,                           Christopher W.%Halladay%NULL%1,                           Emily%O''Neill%NULL%1,                           Stefan%Gravenstein%NULL%1,                           James L.%Rudolph%NULL%1]'
$ws.Range("E28").Value = '[Yongchao%Gao%NULL%1,                           Tao%Liu%NULL%2,                           Weijun%Zhong%NULL%1,                           Rong%Liu%NULL%1,                           Honghao%Zhou%NULL%1,                           Weihua%Huang%endeavor34852@aliyun.com%1,                           Wei%Zhang%csuzhangwei@csu.edu.cn%2]'
$ws.Range("E29").Value = '[Mi Kyung%Kim%NULL%1,                           Jae-Han%Jeon%NULL%2,                           Jae-Han%Jeon%NULL%0,                           Sung-Woo%Kim%NULL%2,                           Sung-Woo%Kim%NULL%0,                           Jun Sung%Moon%NULL%1,                           Nan Hee%Cho%NULL%1,                           Eugene%Han%NULL%1,                           Ji Hong%You%NULL%1,                           Ji Yeon%Lee%NULL%4,                           Miri%Hyun%NULL%1,                           Jae Seok%Park%NULL%4,                           Yong Shik%Kwon%NULL%2,                           Yeon-Kyung%Choi%NULL%1,                           Ki Tae%Kwon%NULL%1,                           Shin Yup%Lee%NULL%1,                           Eon Ju%Jeon%NULL%1,                           Jin-Woo%Kim%NULL%1,                           Hyo-Lim%Hong%NULL%1,                           Hyun Hee%Kwon%NULL%1,                           Chi Young%Jung%NULL%1,                           Yin Young%Lee%NULL%1,                           Eunyeoung%Ha%NULL%1,                           Seung Min%Chung%NULL%1,                           Jian%Hur%NULL%2,                           June Hong%Ahn%NULL%3,                           Na-young%Kim%NULL%1,                           Shin-Woo%Kim%NULL%3,                           Hyun Ha%Chang%NULL%1,                           Yong Hoon%Lee%NULL%1,                           Jaehee%Lee%NULL%1,                           Keun-Gyu%Park%NULL%1,                           Hyun Ah%Kim%NULL%1,                           Ji-Hyun%Lee%NULL%2,                           Ji-Hyun%Lee%NULL%0]'
$ws.Range("E30").Value = '[Pan%Luo%NULL%1,                           Lin%Qiu%NULL%0,                           Yi%Liu%NULL%2,                           Xiu-lan%Liu%NULL%1,                           Jian-ling%Zheng%NULL%1,                           Hui-ying%Xue%NULL%1,                           Wen-hua%Liu%NULL%1,                           Dong%Liu%NULL%3,                           Juan%Li%NULL%3]'
$ws.Range("E31").Value = '[ Z.%Philipose%null%1,                            N.% Smati%null%1,                            C. S. J.% Wong%null%1,                            K.% Aspey%null%1,                            M. A. % Mendall%null%1,         Z.%Philipose%null%1,         N.% Smati%null%1,         C. S. J.% Wong%null%1,         K.% Aspey%null%1,         M. A. % Mendall%null%1]'
$ws.Range("E32").Value = '[ B.%Abu-Jamous%null%1,                            A.% Anisimovich%null%1,                            J.% Baxter%null%1,                            L.% Mackillop%null%1,                            M. P.% Vizcaychipi%null%1,                            A.% McCarthy%null%1,                            R. T. % Khan%null%1,         B.%Abu-Jamous%null%1,         A.% Anisimovich%null%1,         J.% Baxter%null%1,         L.% Mackillop%null%1,         M. P.% Vizcaychipi%null%1,         A.% McCarthy%null%1,         R. T. % Khan%null%1]'
$ws.Range("E33").Value = '[ H. T.%Dashti%null%1,                            D.% Bates%null%1,                            J. M.% Fiskio%null%1,                            E. C.% Roche%null%1,                            S.% Mora%null%1,                            O. % Demler%null%1,         H. T.%Dashti%null%1,         D.% Bates%null%1,         J. M.% Fiskio%null%1,         E. C.% Roche%null%1,         S.% Mora%null%1,         O. % Demler%null%1]'
